# Daily attendance processing - 2026-01-12 18:45:38
# In the "Recorded By" column, swap the order of the two recorders so
# "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$first = $ws.Cells.Find($target)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    $count = 0
    do {
        $current.Value = $replacement
        $count = $count + 1
        $current = $ws.Cells.FindNext($current)
    } while (($current -ne $null) -and ($current.Address() -ne $firstAddress) -and ($count -lt 1000))
}
